# Aggiornamento dati fino al 20/09/2021: aggiunge le righe 375-385 al foglio.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dati da aggiungere: riga foglio, data (seriale Excel), B, C, D
$rows = @(
    @(375, 44449, 0, 0, 0),
    @(376, 44450, 0, 0, 0),
    @(377, 44451, 0, 0, 0),
    @(378, 44452, 0, 0, 0),
    @(379, 44453, 1, 1, 21.81025081788441),
    @(380, 44454, 0, 1, 21.81025081788441),
    @(381, 44455, 0, 1, 21.81025081788441),
    @(382, 44456, 0, 1, 21.81025081788441),
    @(383, 44457, 0, 1, 21.81025081788441),
    @(384, 44458, 0, 1, 21.81025081788441),
    @(385, 44459, 0, 1, 21.81025081788441)
)

# Prende lo stile della cella A374 (colonna data) come riferimento per le nuove celle A.
$dateStyleSource = $ws.Range("A374")

foreach ($r in $rows) {
    $rowIndex = $r[0]

    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = $r[1]

    # Copia solo la formattazione (xlPasteFormats) della cella data di riferimento,
    # cosi' il valore appena scritto non viene sovrascritto.
    $dateStyleSource.Copy()
    $cellA.PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
}

$excel.CutCopyMode = $false
